$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (price and 1h volume % change) per diff
$data = @(
    @{Row=2; Col="D"; Val="306.61"},
    @{Row=2; Col="E"; Val="1.51%"},
    @{Row=3; Col="D"; Val="35.84"},
    @{Row=3; Col="E"; Val="0.97%"},
    @{Row=4; Col="D"; Val="5.065"},
    @{Row=4; Col="E"; Val="0.39%"},
    @{Row=5; Col="D"; Val="0.08086"},
    @{Row=5; Col="E"; Val="0.81%"},
    @{Row=6; Col="D"; Val="1.952"},
    @{Row=6; Col="E"; Val="0.89%"},
    @{Row=7; Col="D"; Val="4.153"},
    @{Row=7; Col="E"; Val="2.23%"},
    @{Row=8; Col="D"; Val="7.779"},
    @{Row=8; Col="E"; Val="-0.40%"},
    @{Row=9; Col="D"; Val="0.9295"},
    @{Row=9; Col="E"; Val="0.47%"},
    @{Row=10; Col="D"; Val="0.1359"},
    @{Row=10; Col="E"; Val="3.59%"},
    @{Row=11; Col="D"; Val="0.1900"},
    @{Row=11; Col="E"; Val="2.40%"},
    @{Row=12; Col="D"; Val="0.09254"},
    @{Row=12; Col="E"; Val="-0.22%"},
    @{Row=13; Col="D"; Val="0.03530"},
    @{Row=13; Col="E"; Val="4.03%"},
    @{Row=14; Col="D"; Val="0.09874"},
    @{Row=14; Col="E"; Val="0.02%"},
    @{Row=15; Col="E"; Val="3.06%"},
    @{Row=16; Col="D"; Val="0.005815"},
    @{Row=16; Col="E"; Val="0.30%"},
    @{Row=17; Col="D"; Val="3.573"},
    @{Row=17; Col="E"; Val="1.77%"},
    @{Row=18; Col="E"; Val="0.58%"},
    @{Row=19; Col="D"; Val="0.3446"},
    @{Row=19; Col="E"; Val="1.27%"},
    @{Row=20; Col="D"; Val="0.1346"},
    @{Row=20; Col="E"; Val="3.40%"},
    @{Row=21; Col="D"; Val="4.905"},
    @{Row=21; Col="E"; Val="-2.80%"},
    @{Row=22; Col="D"; Val="0.2516"},
    @{Row=22; Col="E"; Val="4.76%"},
    @{Row=23; Col="D"; Val="0.04393"},
    @{Row=23; Col="E"; Val="-2.27%"},
    @{Row=24; Col="E"; Val="0.62%"},
    @{Row=25; Col="D"; Val="0.004771"},
    @{Row=26; Col="E"; Val="31.80%"},
    @{Row=27; Col="D"; Val="0.0003130"},
    @{Row=27; Col="E"; Val="4.26%"},
    @{Row=39; Col="D"; Val="0.01969"},
    @{Row=39; Col="E"; Val="2.90%"},
    @{Row=40; Col="D"; Val="0.05009"},
    @{Row=40; Col="E"; Val="5.80%"},
    @{Row=41; Col="D"; Val="0.01111"},
    @{Row=41; Col="E"; Val="15.16%"},
    @{Row=42; Col="D"; Val="0.007628"},
    @{Row=42; Col="E"; Val="3.76%"},
    @{Row=43; Col="D"; Val="0.1378"},
    @{Row=43; Col="E"; Val="3.22%"},
    @{Row=44; Col="D"; Val="0.002102"},
    @{Row=45; Col="D"; Val="0.01081"},
    @{Row=45; Col="E"; Val="-0.68%"},
    @{Row=46; Col="D"; Val="0.00006414"},
    @{Row=46; Col="E"; Val="1.56%"},
    @{Row=47; Col="E"; Val="-0.02%"},
    @{Row=48; Col="D"; Val="63.57"},
    @{Row=48; Col="E"; Val="-1.41%"},
    @{Row=49; Col="D"; Val="0.001191"},
    @{Row=49; Col="E"; Val="-20.02%"},
    @{Row=50; Col="D"; Val="0.00002101"},
    @{Row=50; Col="E"; Val="-0.02%"},
    @{Row=51; Col="D"; Val="0.0002001"},
    @{Row=51; Col="E"; Val="-0.02%"},
)

foreach ($item in $data) {
    $addr = "$($item.Col)$($item.Row)"
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Val
    $cell.Style = "Normal"
}
